# Update the yearly income-statement figures (rows 11-27, columns D:H) on
# the "Overview" sheet with the latest reported period data — mirrors the
# source data refresh described in the commit message ("add monte_carlo
# and update database").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$cols = @("D", "E", "F", "G", "H")

# row -> 5 values for columns D..H (one row per income-statement line item)
$rows = @{
    11 = @(1032970, 1494420, 1782149, 3539419, 6300588)      # فروش
    12 = @(-795466, -1179322, -1391471, -2416864, -4040455)  # بهای تمام شده کالای فروش رفته
    13 = @(237504, 315098, 390678, 1122555, 2260133)         # سود (زیان) ناخالص
    14 = @(-149423, -249510, -188067, -383389, -860439)      # هزینه های عمومی, اداری و تشکیلاتی
    16 = @(4862, -132734, -41754, -13914, -58969)             # خالص سایر درامدها (هزینه ها) ی عملیاتی
    17 = @(92943, -67146, 67378, 644969, 1340725)             # سود (زیان) عملیاتی
    18 = @(-93908, -100098, -78045, -163422, -40103)          # هزینه های مالی
    19 = @(3420, 36660, 49793, -6812, 156252)                 # خالص سایر درامدها و هزینه های غیرعملیاتی
    20 = @(2455, -130584, 39126, 474735, 1456874)             # سود خالص عملیات قبل از مالیات
    21 = @(0, 0, -15968, -153364, -297855)                    # مالیات
    22 = @(2455, -130584, 23158, 321371, 1159019)             # سود خالص عملیات در حال تداوم
    23 = @(0, 0, 0, 0, 0)                                     # سود عملیات متوقف شده
    24 = @(2455, -130584, 23158, 321371, 1159019)             # سود (زیان) خالص
    25 = @(6, -317, 19, 268, 966)                             # سود هر سهم پس از کسر مالیات
    26 = @(412500, 412500, 1200000, 1200000, 1200000)         # سرمایه
    27 = @(2, -109, 19, 268, 966)                             # سود هر سهم بر اساس آخرین سرمایه
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $addr = $cols[$i] + $r
        $ws.Range($addr).Value = $vals[$i]
    }
}

# Row 15 - هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) — D15 keeps its
# literal "-" placeholder; only E15/F15/G15/H15 become numeric.
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = -93479
$ws.Range("G15").Value = -80283
$ws.Range("H15").Value = 0
